# Add the missing Time Recording Log entry (row 18) to the "작성자명" sheet,
# widen column F so the new (longer) activity text fits, and move the
# active-cell selection to the newly filled-in F18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 18: copy the formatting from the row above (row 17), which already
# carries the correct number formats / fonts / borders for each column, then
# overwrite with the new log entry's data. ---
$ws.Range("A17:F17").Copy()
$ws.Range("A18:F18").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A18").Value = 43775                    # 2019-11-06
$ws.Range("B18").Value = 0.79166666666666663      # 19:00
$ws.Range("C18").Value = 0.875                     # 21:00
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 120

$part1 = "Key Class Design, Table Design, Coding Guideline, 추천로직"
$part2 = " 설명 4개 작업물 Review"
$f18 = $ws.Range("F18")
$f18.Value = $part1 + $part2

# Give the trailing part of the text its own run with a distinct font, same
# as the other entries in this log (mixed-font activity descriptions).
$len1 = $part1.Length
$total = ($part1 + $part2).Length
$chars = $f18.Characters($len1 + 1, $total - $len1)
$chars.Font.ColorIndex = -4105   # xlColorIndexAutomatic - keep default color
$chars.Font.Name = "맑은 고딕"
$chars.Font.Size = 10

# --- Widen column F to fit the new, longer activity text. ---
$ws.Columns.Item(6).ColumnWidth = 74

# --- Move the selection to the newly entered cell. ---
$ws.Activate() | Out-Null
$ws.Range("F18").Select() | Out-Null
